$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns("D:D").Insert()
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$ws.Range("D1").ColumnWidth = $ws.Range("E1").ColumnWidth
